$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header row
$ws.Range("B1").Value = "Sucursal"
$ws.Range("C1").Value = "CARPETA DE GESTIÓN ELECTRO"
$ws.Range("D1").Value = "Comentarios Carpeta Electro"
$ws.Range("E1").Value = "Pregunta"

# Row 2
$ws.Range("B2").Value = "La paz"
$ws.Range("C2").Value = "¿Tiene firmados los objetivos de todos los vendedores?;¿Tiene planificación de trabajo por el desvío de objetivos -mes anterior? (template);¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;¿Tiene cierre y devoluciones realizadas mes anterior?;Comunicación: firma de procesos claves;"
$ws.Range("D2").Value = "El chevk list de esta semana falta por encontrarse de vacaciones"
$ws.Range("E2").ClearContents()

# Row 3
$ws.Range("B3").Value = "La paz"
$ws.Range("C3").Value = "¿Tiene firmados los objetivos de todos los vendedores?;¿Tiene planificación de trabajo por el desvío de objetivos -mes anterior? (template);¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;¿Tiene cierre y devoluciones realizadas mes anterior?;Comunicación: firma de procesos claves;Tiene el Gerente realizado el check list del día?;"
$ws.Range("D3").Value = "Se realizó la reunión con los vendedores haciendo hincapié en blister y garantía ."
$ws.Range("E3").ClearContents()

# Row 4
$ws.Range("B4").Value = "Formosa "
$ws.Range("C4").Value = "¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;Comunicación: firma de procesos claves;Tiene el Gerente realizado el check list del día?;"
$ws.Range("D4").ClearContents()
